$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new display text, exactly as the target diff specifies.
# All of these columns (D = Price, E = Volume(1h)) store plain text in the
# workbook (numbers/percentages kept as literal strings), and columns B/C
# hold the coin name / link text for the two swapped rows (41/42).
$updates = [ordered]@{
    "D2" = "'255.81"
    "E2" = "'4.24%"
    "D3" = "'27.57"
    "E3" = "'-2.59%"
    "E4" = "'-0.71%"
    "D5" = "'0.05914"
    "E5" = "'3.79%"
    "D6" = "'6.686"
    "E6" = "'0.81%"
    "D7" = "'0.8653"
    "E7" = "'1.73%"
    "D8" = "'1.029"
    "E8" = "'13.49%"
    "D9" = "'0.1419"
    "E9" = "'3.68%"
    "E10" = "'8.38%"
    "D11" = "'0.07206"
    "E11" = "'1.94%"
    "D12" = "'0.03257"
    "E12" = "'2.21%"
    "D13" = "'0.09220"
    "E13" = "'-0.03%"
    "D14" = "'0.001554"
    "E14" = "'2.07%"
    "D15" = "'0.0006057"
    "E15" = "'-93.95%"
    "D16" = "'0.005783"
    "E16" = "'-2.43%"
    "D17" = "'3.484"
    "E17" = "'-0.21%"
    "D18" = "'3.267"
    "E18" = "'1.98%"
    "E21" = "'2.49%"
    "D22" = "'3.533"
    "E22" = "'0.32%"
    "D23" = "'0.04171"
    "E23" = "'2.41%"
    "E24" = "'1.59%"
    "D25" = "'0.001221"
    "E25" = "'-0.01%"
    "E26" = "'8.81%"
    "D27" = "'0.0001202"
    "E27" = "'0.25%"
    "E28" = "'33.96%"
    "D40" = "'0.03812"
    "E40" = "'-0.36%"
    "B41" = "KickToken"
    "C41" = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
    "D41" = "'0.005511"
    "E41" = "'7.49%"
    "B42" = "BKEXToken"
    "C42" = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
    "D42" = "'0.1101"
    "E42" = "'3.33%"
    "D43" = "'0.001904"
    "E43" = "'-13.43%"
    "D44" = "'0.01072"
    "E44" = "'16.66%"
    "D45" = "'0.00005439"
    "E45" = "'3.11%"
    "E46" = "'0.23%"
    "E47" = "'4.06%"
    "D48" = "'0.002164"
    "E48" = "'-4.62%"
    "D49" = "'0.00002104"
    "E49" = "'0.23%"
    "D50" = "'0.0002004"
    "E50" = "'0.23%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = $updates[$addr]
    # Re-apply the default style: entering a leading-apostrophe literal for
    # numeric-looking text (prices / percentages) marks the cell with Excel's
    # "quote prefix" flag; these cells were plain, unstyled text before the
    # edit, so restore that so only the cell *content* changes, matching the diff.
    $cell.Style = "Normal"
}
